# ============================================================
# Insert a new "2022-Q3" quarter into the 长城汽车 (02333) workbook.
#
#  1. "总计" (summary) sheet gets a brand-new row 2 with the 2022-Q3
#     totals; every later row keeps its data but shifts down one slot.
#  2. A brand-new worksheet named "2022-Q3" is inserted right after
#     "总计" holding the per-fund holdings detail for that quarter
#     (same shape as the existing 2022-Q2 / 2022-Q1 / ... sheets).
# ============================================================

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)   # "总计"

# ---------- 1. "总计" sheet: add the 2022-Q3 row -----------------------
# Target final table (A=index, B=quarter, C=count, D=market value):
#   0  2022-Q3  23  11.99
#   1  2022-Q2  56  47.98
#   2  2022-Q1  36  29.92
#   3  2021-Q4  73  113.85
#   4  2021-Q3  80  113.33
#   5  2021-Q2  55  56.48
#   6  2021-Q1  36  41.06
#   7  2020-Q4  86  51.77
# Row 9 is brand new, so first clone row 8's formatting down into it,
# then overwrite every data row (2-9) with its final values.

$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122)   # xlPasteFormats

$summaryData = @(
    @(0, "2022-Q3", 23, 11.99),
    @(1, "2022-Q2", 56, 47.98),
    @(2, "2022-Q1", 36, 29.92),
    @(3, "2021-Q4", 73, 113.85),
    @(4, "2021-Q3", 80, 113.33),
    @(5, "2021-Q2", 55, 56.48),
    @(6, "2021-Q1", 36, 41.06),
    @(7, "2020-Q4", 86, 51.77)
)

$r = 2
foreach ($row in $summaryData) {
    $summary.Cells.Item($r, 1).Value = $row[0]
    $summary.Cells.Item($r, 2).Value = $row[1]
    $summary.Cells.Item($r, 3).Value = $row[2]
    $summary.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# ---------- 2. Brand-new "2022-Q3" detail sheet -------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Reuse the bold/bordered header style already used throughout the
# workbook (copied from the "总计" header) for the header row...
$summary.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)    # xlPasteFormats

# ...and for the index column (A2:A24).
$summary.Range("A2").Copy()
$q3.Range("A2:A24").PasteSpecial(-4122)   # xlPasteFormats

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3Data = @(
    @(0, "010902", "博时成长领航灵活配置混合A", "46.83", "81.70", "5.28", "2.4726", 3),
    @(1, "012344", "嘉实领先优势混合A", "46.88", "87.50", "4.69", "2.1987", 6),
    @(2, "010041", "嘉实港股优势混合A", "40.43", "89.69", "4.61", "1.8638", 4),
    @(3, "160527", "博时研究优选3年封闭运作灵活配置混合A", "16.93", "97.11", "7.64", "1.2935", 3),
    @(4, "001878", "嘉实沪港深精选股票", "22.02", "91.33", "4.71", "1.0371", 5),
    @(5, "160726", "嘉实瑞享定期开放灵活配置混合", "12.64", "83.43", "5.95", "0.7521", 2),
    @(6, "009138", "嘉实瑞成两年持有期混合A", "10.97", "90.59", "5.91", "0.6483", 2),
    @(7, "010903", "博时成长领航灵活配置混合C", "6.84", "81.70", "5.28", "0.3612", 3),
    @(8, "010042", "嘉实港股优势混合C", "5.20", "89.69", "4.61", "0.2397", 4),
    @(9, "013123", "汇添富精选核心优势一年持有混合A", "5.35", "80.40", "4.36", "0.2333", 4),
    @(10, "005228", "汇添富港股通专注成长混合", "4.41", "74.46", "4.56", "0.2011", 6),
    @(11, "009139", "嘉实瑞成两年持有期混合C", "2.99", "90.59", "5.91", "0.1767", 2),
    @(12, "006511", "博道卓远混合A", "3.37", "80.04", "3.51", "0.1183", 5),
    @(13, "014506", "博时成长臻选混合A", "3.12", "86.00", "3.37", "0.1051", 9),
    @(14, "013550", "汇添富品牌价值一年持有混合A", "2.12", "78.71", "3.81", "0.0808", 4),
    @(15, "011924", "嘉实港股互联网产业核心资产混合A", "1.20", "87.88", "5.58", "0.0670", 4),
    @(16, "160528", "博时研究优选3年封闭运作灵活配置混合C", "0.59", "97.11", "7.64", "0.0451", 3),
    @(17, "012345", "嘉实领先优势混合C", "0.69", "87.50", "4.69", "0.0324", 6),
    @(18, "011925", "嘉实港股互联网产业核心资产混合C", "0.41", "87.88", "5.58", "0.0229", 4),
    @(19, "013124", "汇添富精选核心优势一年持有混合C", "0.27", "80.40", "4.36", "0.0118", 4),
    @(20, "013551", "汇添富品牌价值一年持有混合C", "0.28", "78.71", "3.81", "0.0107", 4),
    @(21, "006512", "博道卓远混合C", "0.28", "80.04", "3.51", "0.0098", 5),
    @(22, "014507", "博时成长臻选混合C", "0.14", "86.00", "3.37", "0.0047", 9)
)

$r = 2
foreach ($row in $q3Data) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = "'" + $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = "'" + $row[3]
    $q3.Cells.Item($r, 5).Value = "'" + $row[4]
    $q3.Cells.Item($r, 6).Value = "'" + $row[5]
    $q3.Cells.Item($r, 7).Value = "'" + $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r++
}
